$wb = $excel.ActiveWorkbook

# --- Sheet "Description": update generated timestamp + summary stats ---
$wsDesc = $wb.Worksheets.Item("Description")
$wsDesc.Range("B2").Value = "2025-10-09 11:37:07"
$wsDesc.Range("B13").Value = 3690
$wsDesc.Range("B20").Value = 11

# --- Sheet "Overlap Analysis": refresh overlap-analysis dataset (rows 2:12) ---
$wsOv = $wb.Worksheets.Item("Overlap Analysis")

$overlapRows = @(
    @{Row=2; A="52060800020002"; B="52060800021005"; C="41224.30805170261"; D="42167.66026968733"; E="10100.23377610135"; F="0.2450067509546517"; G="0.2395255916857689"}
    @{Row=3; A="52060800020003"; B="52060800022008"; C="32303.77980558742"; D="145182.3518493054"; E="15434.23423791282"; F="0.4777841581016237"; G="0.1063093002786802"}
    @{Row=4; A="52060800050012"; B="52060800050011"; C="122400.297211675"; D="270263.5685508295"; E="122400.297211675"; F="1"; G="0.4528923297653223"}
    @{Row=5; A="52060800050017"; B="52060800050016"; C="15423.42325961922"; D="2883.437007030753"; E="2791.950138016949"; F="0.1810201335346015"; G="0.9682715908859012"}
    @{Row=6; A="52060800050018"; B="52060800050016"; C="7540.410064111213"; D="1186.639068642798"; E="1186.639068642798"; F="0.1573706281957581"; G="1"}
    @{Row=7; A="52060800050018"; B="52060800050015"; C="7540.410064111213"; D="0.02071007287384314"; E="0.02071007287384314"; F="2.746544643826905e-06"; G="1"}
    @{Row=8; A="52060800050006"; B="52060800051004"; C="17011.19436189108"; D="1672752.195255189"; E="5709.037937349282"; F="0.3356047680072846"; G="0.003412960959515185"}
    @{Row=9; A="52060800050023"; B="52060800050005"; C="44954.81772431835"; D="105421.6830197758"; E="27763.44943372393"; F="0.6175856301760794"; G="0.2633561582252094"}
    @{Row=10; A="52060800050005"; B="52060800051004"; C="97407.99722729861"; D="1672752.195255189"; E="28268.63710708431"; F="0.2902085856577084"; G="0.0168994769143147"}
    @{Row=11; A="52060800080003"; B="52060800081001"; C="10663.37631213075"; D="139176.3641311153"; E="3128.28280175185"; F="0.2933670078015618"; G="0.02247711255630128"}
    @{Row=12; A="52060800120002"; B="52060800120003"; C="14675.98062515429"; D="26449.51773102354"; E="3853.087095435435"; F="0.2625437573030952"; G="0.1456770265007901"}
)

foreach ($rowData in $overlapRows) {
    $r = $rowData.Row
    # A/B hold zero-padded-looking numeric ID strings; route the literal
    # text through a TEXT() formula + paste-as-values so Excel keeps them
    # as text instead of auto-coercing to numbers.
    $wsOv.Cells.Item($r, 1).Formula = "=TEXT(" + $rowData.A + ",""0"")"
    $wsOv.Cells.Item($r, 2).Formula = "=TEXT(" + $rowData.B + ",""0"")"
    $wsOv.Cells.Item($r, 3).Value = [double]$rowData.C
    $wsOv.Cells.Item($r, 4).Value = [double]$rowData.D
    $wsOv.Cells.Item($r, 5).Value = [double]$rowData.E
    $wsOv.Cells.Item($r, 6).Value = [double]$rowData.F
    $wsOv.Cells.Item($r, 7).Value = [double]$rowData.G
}

$wsOv.Range("A2:B12").Copy()
$wsOv.Range("A2:B12").PasteSpecial(-4163)
$excel.CutCopyMode = $false

